# Update the generated three-digit x one-digit multiplication answers
# in the table to match the newly generated set of problems/answers.

$d = $word.ActiveDocument

$replacements = @(
    @("425×7=2975", "283×8=2264"),
    @("213×7=1491", "496×9=4464"),
    @("555×5=2775", "568×9=5112"),
    @("956×2=1912", "851×6=5106"),
    @("940×5=4700", "996×5=4980"),
    @("766×4=3064", "843×9=7587"),
    @("108×7=756",  "371×9=3339"),
    @("286×8=2288", "895×8=7160"),
    @("148×9=1332", "693×2=1386"),
    @("651×2=1302", "297×7=2079"),
    @("504×7=3528", "942×3=2826"),
    @("472×2=944",  "796×6=4776"),
    @("988×3=2964", "854×6=5124"),
    @("131×5=655",  "446×3=1338"),
    @("521×6=3126", "736×8=5888"),
    @("798×6=4788", "227×7=1589"),
    @("988×6=5928", "282×5=1410"),
    @("625×3=1875", "708×2=1416"),
    @("186×6=1116", "654×3=1962"),
    @("404×4=1616", "923×8=7384"),
    @("780×7=5460", "922×3=2766"),
    @("903×9=8127", "357×2=714"),
    @("838×9=7542", "222×6=1332"),
    @("259×3=777",  "481×4=1924"),
    @("287×4=1148", "521×7=3647")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
